$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'314.44"
$ws.Range("E2").Value = "'1.94%"
# Row 3
$ws.Range("D3").Value = "'39.28"
$ws.Range("E3").Value = "'-1.58%"
# Row 4
$ws.Range("D4").Value = "'5.151"
$ws.Range("E4").Value = "'-0.07%"
# Row 5
$ws.Range("D5").Value = "'0.08170"
$ws.Range("E5").Value = "'0.35%"
# Row 6
$ws.Range("D6").Value = "'1.985"
$ws.Range("E6").Value = "'1.94%"
# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.328"
$ws.Range("E7").Value = "'2.11%"
# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9386"
$ws.Range("E8").Value = "'0.96%"
# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1310"
$ws.Range("E9").Value = "'-8.87%"
# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1973"
$ws.Range("E10").Value = "'2.65%"
# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09024"
$ws.Range("E11").Value = "'-0.79%"
# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03527"
$ws.Range("E12").Value = "'0.33%"
# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09744"
$ws.Range("E13").Value = "'-0.44%"
# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001406"
$ws.Range("E14").Value = "'0.98%"
# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006048"
$ws.Range("E15").Value = "'3.30%"
# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.632"
$ws.Range("E16").Value = "'-7.34%"
# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.381"
$ws.Range("E17").Value = "'3.33%"
# Row 18
$ws.Range("D18").Value = "'3.117"
$ws.Range("E18").Value = "'-6.25%"
# Row 19
$ws.Range("D19").Value = "'0.3470"
$ws.Range("E19").Value = "'1.16%"
# Row 20
$ws.Range("D20").Value = "'0.1316"
$ws.Range("E20").Value = "'0.26%"
# Row 21
$ws.Range("E21").Value = "'7.31%"
# Row 22
$ws.Range("E22").Value = "'2.61%"
# Row 23
$ws.Range("D23").Value = "'0.04380"
$ws.Range("E23").Value = "'0.06%"
# Row 24
$ws.Range("D24").Value = "'0.001240"
$ws.Range("E24").Value = "'0.89%"
# Row 25
$ws.Range("D25").Value = "'0.004764"
$ws.Range("E25").Value = "'8.98%"
# Row 26
$ws.Range("D26").Value = "'0.0003896"
$ws.Range("E26").Value = "'199.69%"
# Row 27
$ws.Range("E27").Value = "'-7.65%"
# Row 39
$ws.Range("D39").Value = "'0.02246"
$ws.Range("E39").Value = "'9.51%"
# Row 40
$ws.Range("D40").Value = "'0.05197"
# Row 41
$ws.Range("D41").Value = "'0.007741"
$ws.Range("E41").Value = "'4.61%"
# Row 42
$ws.Range("D42").Value = "'0.01031"
$ws.Range("E42").Value = "'4.50%"
# Row 43
$ws.Range("D43").Value = "'0.1398"
$ws.Range("E43").Value = "'2.41%"
# Row 44
$ws.Range("D44").Value = "'0.002103"
$ws.Range("E44").Value = "'-1.25%"
# Row 45
$ws.Range("D45").Value = "'0.008872"
$ws.Range("E45").Value = "'-5.30%"
# Row 46
$ws.Range("D46").Value = "'0.00006822"
$ws.Range("E46").Value = "'7.11%"
# Row 47
$ws.Range("E47").Value = "'0.06%"
# Row 48
$ws.Range("D48").Value = "'0.003010"
$ws.Range("E48").Value = "'10.90%"
# Row 49
$ws.Range("E49").Value = "'30.04%"
# Row 50
$ws.Range("E50").Value = "'0.06%"
# Row 51
$ws.Range("E51").Value = "'0.06%"
